$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction8")

$values = @(2,3,5,8,9,11,12,13,14,16,17,19,20,21,22,23,25,26,28,30,31,32,33,34,35,36,38,39,41,42,43,44,45,46,47,50,51,53,54,55,56,57,58,59,61,65,66,68,69,70,72,73)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
